$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: new entry for 2023-02-09 (date serial 44966), afternoon AI-agent work session ---
$ws.Range("A11").Value = 44966
$ws.Range("B11").Value = "13.30 - 16.00"
$ws.Range("D11").Value = 2.5
$ws.Range("E11").Value = "Working on AI agents"

# Re-use the existing date number format from the column above for the new date cells
$dateFmt = $ws.Range("A11").NumberFormat

# --- Row 12: entry for 2023-02-10 (date serial 44967), morning AI-navigation session ---
$ws.Range("A12").Value = 44967
$ws.Range("A12").NumberFormat = $dateFmt
$ws.Range("E12").Value = "Implemented basic AI navigation."
$ws.Range("B12").Value = "9.30 - 12.30"
$ws.Range("D12").Value = 3

# --- Row 13: entry for 2023-02-10 (date serial 44967), afternoon flocking/boids research ---
$ws.Range("A13").Value = 44967
$ws.Range("A13").NumberFormat = $dateFmt
$ws.Range("B13").Value = "15.00 - 17.00"
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = "Researching flocking/boids"

# Move the cell cursor to where the author last left it
$null = $ws.Range("AF33").Select()
